$wb = $excel.ActiveWorkbook

# --- orb_mission sheet (data + formula + hyperlink updates) ---
$ws = $wb.Worksheets.Item("orb_mission")

# Row 5: mom dump related slew time -> 100
$ws.Range("B5").Value = 100

# Row 6: slew_angle -> 135
$ws.Range("B6").Value = 135

# Row 7: slew_time -> 1200 sec, add minutes formula + label + old value reference
$ws.Range("B7").Value = 1200
$ws.Range("D7").Formula = "=B7/60"
$ws.Range("E7").Value = "min"
$ws.Range("F7").Value = 3000

# Row 10: app_slew -> 150, plus hyperlink reference in H10
$ws.Range("B10").Value = 150

# Row 11: app_time -> 600
$ws.Range("B11").Value = 600

# Add hyperlink (source citation) on H10 - add after the "min" shared string
# above so the shared-string table order matches (min, then URL).
$null = $ws.Hyperlinks.Add($ws.Range("H10"), "http://www.msss.com/mars/global_surveyor/mgs_msn_plan/section5/section5.html")

# Make orb_mission the active sheet/tab, with B6 selected
$ws.Activate()
$null = $ws.Range("B6").Select()
